$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting, since many values
# (e.g. '61.298.02') use '.' as thousands separators and are not valid numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range('D2').Value = '61.298.02'
$ws.Range('E2').Value = '  -1.65%  '

# Row 3
$ws.Range('D3').Value = '2.981.09'
$ws.Range('E3').Value = '  -1.29%  '

# Row 4
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.18%  '

# Row 5
$ws.Range('D5').Value = '594.93'
$ws.Range('E5').Value = '  +1.73%  '

# Row 6
$ws.Range('D6').Value = '143.52'
$ws.Range('E6').Value = '  -2.68%  '

# Row 7
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.16%  '

# Row 8
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').Value = '0.513'
$ws.Range('E8').Value = '  -2.12%  '

# Row 9
$ws.Range('B9').Value = 'LidoStakedEther'
$ws.Range('C9').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D9').Value = '2.979.65'
$ws.Range('E9').Value = '  -1.19%  '

# Row 10
$ws.Range('D10').Value = '0.148'
$ws.Range('E10').Value = '  -0.67%  '

# Row 11
$ws.Range('D11').Value = '6.05'
$ws.Range('E11').Value = '  +3.78%  '

# Row 12
$ws.Range('D12').Value = '0.452'
$ws.Range('E12').Value = '  +1.74%  '

# Row 13
$ws.Range('D13').Value = '0.0000227'
$ws.Range('E13').Value = '  -0.79%  '

# Row 14
$ws.Range('D14').Value = '34.07'
$ws.Range('E14').Value = '  -1.57%  '

# Row 15
$ws.Range('E15').Value = '  +3.18%  '

# Row 16
$ws.Range('D16').Value = '3.469.40'
$ws.Range('E16').Value = '  -1.41%  '

# Row 17
$ws.Range('D17').Value = '6.90'
$ws.Range('E17').Value = '  -2.46%  '

# Row 18
$ws.Range('D18').Value = '61.201.33'
$ws.Range('E18').Value = '  -1.77%  '

# Row 19
$ws.Range('D19').Value = '2.976.22'
$ws.Range('E19').Value = '  -1.52%  '

# Row 20
$ws.Range('D20').Value = '447.75'
$ws.Range('E20').Value = '  -3.40%  '

# Row 21
$ws.Range('D21').Value = '13.95'
$ws.Range('E21').Value = '  -0.26%  '

# Row 22
$ws.Range('D22').Value = '0.680'
$ws.Range('E22').Value = '  -0.63%  '

# Row 23
$ws.Range('D23').Value = '7.33'
$ws.Range('E23').Value = '  -1.71%  '

# Row 24
$ws.Range('D24').Value = '81.34'
$ws.Range('E24').Value = '  -0.51%  '

# Row 25
$ws.Range('D25').Value = '10.61'
$ws.Range('E25').Value = '  +6.20%  '

# Row 26
$ws.Range('D26').Value = '2.18'
$ws.Range('E26').Value = '  -3.37%  '

# Row 27
$ws.Range('D27').Value = '11.98'
$ws.Range('E27').Value = '  -2.55%  '

# Row 28
$ws.Range('E28').Value = '  +0.11%  '

# Row 29
$ws.Range('D29').Value = '2.68'
$ws.Range('E29').Value = '  +1.99%  '

# Row 30
$ws.Range('E30').Value = '  -0.13%  '

# Row 31
$ws.Range('D31').Value = '7.17'
$ws.Range('E31').Value = '  +0.08%  '

# Row 32
$ws.Range('D32').Value = '2.06'
$ws.Range('E32').Value = '  -1.99%  '

# Row 33
$ws.Range('D33').Value = '27.14'
$ws.Range('E33').Value = '  -7.01%  '

# Row 34
$ws.Range('E34').Value = '  +1.39%  '

# Row 35
$ws.Range('D35').Value = '0.0₃0814'
$ws.Range('E35').Value = '  +1.99%  '

# Row 36
$ws.Range('D36').Value = '1.02'
$ws.Range('E36').Value = '  -1.10%  '

# Row 37
$ws.Range('D37').Value = '5.76'
$ws.Range('E37').Value = '  -0.16%  '

# Row 38
$ws.Range('D38').Value = '50.12'
$ws.Range('E38').Value = '  -0.41%  '

# Row 39
$ws.Range('D39').Value = '8.95'
$ws.Range('E39').Value = '  -1.07%  '

# Row 40
$ws.Range('E40').Value = '  -5.19%  '

# Row 41
$ws.Range('D41').Value = '0.123'
$ws.Range('E41').Value = '  +9.09%  '

# Row 42
$ws.Range('D42').Value = '2.84'
$ws.Range('E42').Value = '  -3.80%  '

# Row 43
$ws.Range('D43').Value = '384.97'
$ws.Range('E43').Value = '  -0.97%  '

# Row 44
$ws.Range('B44').Value = 'Arweave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D44').Value = '39.41'
$ws.Range('E44').Value = '  +5.14%  '

# Row 45
$ws.Range('D45').Value = '0.268'
$ws.Range('E45').Value = '  -2.27%  '

# Row 46
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').Value = '0.0348'
$ws.Range('E46').Value = '  -2.51%  '

# Row 47
$ws.Range('D47').Value = '2.681.23'
$ws.Range('E47').Value = '  -2.73%  '

# Row 48
$ws.Range('D48').Value = '130.21'
$ws.Range('E48').Value = '  +1.45%  '

# Row 49
$ws.Range('E49').Value = '  +0.14%  '

# Row 50
$ws.Range('E50').Value = '  -1.56%  '

# Row 51
$ws.Range('D51').Value = '2.14'
$ws.Range('E51').Value = '  -1.14%  '
